# MainReport.xlsx update: add three new rows of data (rows 8-10) capturing a
# new "L120.csv" class with MDF/RMS/MDF features and Grid1/Grid2 = 32
# settings, matching the author's "major update, added new class" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: L120.csv / MDF -------------------------------------------------
$ws.Cells.Item(8,1).Value  = "L120.csv"
$ws.Cells.Item(8,2).Value  = "MDF"
$ws.Cells.Item(8,3).Value  = 20282
$ws.Cells.Item(8,4).Value  = 250
$ws.Cells.Item(8,5).Value  = "Grid1: 32"
$ws.Cells.Item(8,6).Value  = 4.9997147825811075
$ws.Cells.Item(8,7).Value  = 1.0123546844457523
$ws.Cells.Item(8,8).Value  = 2.1850623651544141
$ws.Cells.Item(8,9).Value  = 2.1877179967923386
$ws.Cells.Item(8,10).Value = 153.13073432744957
$ws.Cells.Item(8,11).Value = 500
$ws.Cells.Item(8,12).Value = 2.4913278861865829
$ws.Cells.Item(8,13).Value = 4.4964457027896536
$ws.Cells.Item(8,14).Value = "Grid2: 32"
$ws.Cells.Item(8,15).Value = 4.957208825138796
$ws.Cells.Item(8,16).Value = 12.455041257342652
$ws.Cells.Item(8,17).Value = 2.116146697895704
$ws.Cells.Item(8,18).Value = 2.1757732740706923
$ws.Cells.Item(8,19).Value = 130.66121667900043
$ws.Cells.Item(8,20).Value = 500
$ws.Cells.Item(8,21).Value = 2.4563775322756314
$ws.Cells.Item(8,22).Value = 4.2481080206193091
$ws.Cells.Item(8,23).Value = 0

# --- Row 9: L120.csv / RMS -------------------------------------------------
$ws.Cells.Item(9,1).Value  = "L120.csv"
$ws.Cells.Item(9,2).Value  = "RMS"
$ws.Cells.Item(9,3).Value  = 20282
$ws.Cells.Item(9,4).Value  = 250
$ws.Cells.Item(9,5).Value  = "Grid1: 32"
$ws.Cells.Item(9,6).Value  = 4.9827541319052324
$ws.Cells.Item(9,7).Value  = 7.8740373620503528
$ws.Cells.Item(9,8).Value  = -0.22747177973472568
$ws.Cells.Item(9,9).Value  = -0.20824049978381717
$ws.Cells.Item(9,10).Value = 0.59228157182295083
$ws.Cells.Item(9,11).Value = 1.65
$ws.Cells.Item(9,12).Value = 2.5041788721098985
$ws.Cells.Item(9,13).Value = 4.3239160808207764
$ws.Cells.Item(9,14).Value = "Grid2: 32"
$ws.Cells.Item(9,15).Value = 4.7404864785606913
$ws.Cells.Item(9,16).Value = 31.350721975264484
$ws.Cells.Item(9,17).Value = 0.056693418498374525
$ws.Cells.Item(9,18).Value = -0.091493144296735884
$ws.Cells.Item(9,19).Value = 1.1394451350814223
$ws.Cells.Item(9,20).Value = 2.25
$ws.Cells.Item(9,21).Value = 2.729627503259529
$ws.Cells.Item(9,22).Value = 4.9528656752405942
$ws.Cells.Item(9,23).Value = 0

# --- Row 10: L120.csv / MDF -------------------------------------------------
$ws.Cells.Item(10,1).Value  = "L120.csv"
$ws.Cells.Item(10,2).Value  = "MDF"
$ws.Cells.Item(10,3).Value  = 20282
$ws.Cells.Item(10,4).Value  = 1000
$ws.Cells.Item(10,5).Value  = "Grid1: 32"
$ws.Cells.Item(10,6).Value  = 4.9899433448642005
$ws.Cells.Item(10,7).Value  = 5.9604783568598751
$ws.Cells.Item(10,8).Value  = 2.0400266820151334
$ws.Cells.Item(10,9).Value  = 2.0344533587457043
$ws.Cells.Item(10,10).Value = 109.65455632127134
$ws.Cells.Item(10,11).Value = 500
$ws.Cells.Item(10,12).Value = 2.5060892817466955
$ws.Cells.Item(10,13).Value = 4.4042009296413376
$ws.Cells.Item(10,14).Value = "Grid2: 32"
$ws.Cells.Item(10,15).Value = 4.9517624780973328
$ws.Cells.Item(10,16).Value = 13.285536622456886
$ws.Cells.Item(10,17).Value = 2.0679762691188244
$ws.Cells.Item(10,18).Value = 2.1299249175009325
$ws.Cells.Item(10,19).Value = 116.94354885435517
$ws.Cells.Item(10,20).Value = 500
$ws.Cells.Item(10,21).Value = 2.4427652631336318
$ws.Cells.Item(10,22).Value = 4.2369687702681196
$ws.Cells.Item(10,23).Value = 0

# --- Column widths (E and N grew slightly to fit "Grid1: 32" text) --------
$ws.Columns.Item(5).ColumnWidth = 7.61
$ws.Columns.Item(14).ColumnWidth = 7.61

# --- View: scrolled right so column I is the left-most visible column -----
$ws.Range("I1").Select()

# --- Selection returns to E9 to match the saved cursor position -----------
$ws.Range("E9").Select()
